# Rename the inline pictures' internal "name" (wp:docPr / pic:cNvPr @name)
# in the document's headers/footers:
#   - header (first-page): BTec_Logo-Orange      image1.jpg -> image2.jpg
#   - footer (default):    PearsonLogo (id=2)    image2.png -> image1.png
#   - footer (first-page): PearsonLogo (id=3)    image2.png -> image1.png
#
# wdHeaderFooterIndex: 1 = wdHeaderFooterPrimary, 2 = wdHeaderFooterFirstPage

$d = $word.ActiveDocument

function Rename-InlineLogo {
    param($range, [string]$newName, [string]$expectedAlt)

    for ($i = 1; $i -le $range.InlineShapes.Count; $i++) {
        $shape = $range.InlineShapes($i)
        if ($shape.AlternativeText -eq $expectedAlt) {
            $shape.Name = $newName
        }
    }
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections($s)

    # Header, first page -> BTec_Logo-Orange: image1.jpg -> image2.jpg
    $header = $section.Headers(2)
    if ($header.Exists) {
        Rename-InlineLogo $header.Range "image2.jpg" "BTec_Logo-Orange"
    }

    # Footer, default (primary) -> PearsonLogo: image2.png -> image1.png
    $footerDefault = $section.Footers(1)
    if ($footerDefault.Exists) {
        Rename-InlineLogo $footerDefault.Range "image1.png" "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png"
    }

    # Footer, first page -> PearsonLogo: image2.png -> image1.png
    $footerFirst = $section.Footers(2)
    if ($footerFirst.Exists) {
        Rename-InlineLogo $footerFirst.Range "image1.png" "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png"
    }
}
